# Weekly update: insert a new price record for Poroto granado at row 36,
# pushing the existing rows (36-67) down by one (to 37-68).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 36; Excel shifts rows 36:67 down to 37:68.
$ws.Rows("36:36").Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Cells.Item(36, 1).Value = 2
$ws.Cells.Item(36, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(36, 3).Value = "Coquimbo"
$ws.Cells.Item(36, 4).Value = 44587
$ws.Cells.Item(36, 5).Value = 4
$ws.Cells.Item(36, 6).Value = 100112030
$ws.Cells.Item(36, 7).Value = "Poroto granado"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 600
$ws.Cells.Item(36, 11).Value = 23000
$ws.Cells.Item(36, 12).Value = 25000
$ws.Cells.Item(36, 13).Value = 24000
$ws.Cells.Item(36, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(36, 16).Value = 960
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
